$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.069299666666667
$ws.Range("H2").Value = 3.207899
$ws.Range("I2").Value = 0.003616700200628781
$ws.Range("J2").Value = 0.003616700200628781
$ws.Range("M2").Value = 1.743137
$ws.Range("N2").Value = 5.229411
$ws.Range("O2").Value = 0.03144673183548247
$ws.Range("P2").Value = 0.03144673183548247
$ws.Range("Q2").Value = 1.863935813054334
$ws.Range("R2").Value = 16.775422317489
$ws.Range("S2").Value = 0.0001137334013385089
$ws.Range("T2").Value = 0.0001137334013385089
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.069299666666667
$ws.Range("H3").Value = 3.207899
$ws.Range("I3").Value = 0.003616700200628781
$ws.Range("J3").Value = 0.003616700200628781
$ws.Range("N3").Value = 3.848628
$ws.Range("O3").Value = 0.02314348071905789
$ws.Range("P3").Value = 0.02314348071905789
$ws.Range("Q3").Value = 1.371778879174667
$ws.Range("R3").Value = 12.346009912572
$ws.Range("S3").Value = 0.000083703031359865
$ws.Range("T3").Value = 0.00008370303135986499
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.069299666666667
$ws.Range("H4").Value = 3.207899
$ws.Range("I4").Value = 0.003616700200628781
$ws.Range("J4").Value = 0.003616700200628781
$ws.Range("M4").Value = 52.405407
$ws.Range("N4").Value = 157.216221
$ws.Range("O4").Value = 0.9454097874454597
$ws.Range("P4").Value = 0.9454097874454597
$ws.Range("Q4").Value = 56.03708423663101
$ws.Range("R4").Value = 504.3337581296791
$ws.Range("S4").Value = 0.003419263767930407
$ws.Range("T4").Value = 0.003419263767930407
$ws.Range("I5").Value = 0.8238194745364892
$ws.Range("J5").Value = 0.8238194745364891
$ws.Range("M5").Value = 1.743137
$ws.Range("N5").Value = 5.229411
$ws.Range("O5").Value = 0.03144673183548247
$ws.Range("P5").Value = 0.03144673183548247
$ws.Range("Q5").Value = 424.5711662285977
$ws.Range("R5").Value = 3821.140496057379
$ws.Range("S5").Value = 0.02590643009659706
$ws.Range("T5").Value = 0.02590643009659705
$ws.Range("I6").Value = 0.8238194745364892
$ws.Range("J6").Value = 0.8238194745364891
$ws.Range("N6").Value = 3.848628
$ws.Range("O6").Value = 0.02314348071905789
$ws.Range("P6").Value = 0.02314348071905789
$ws.Range("R6").Value = 2812.199749658292
$ws.Range("S6").Value = 0.01906605012491964
$ws.Range("T6").Value = 0.01906605012491964
$ws.Range("I7").Value = 0.8238194745364892
$ws.Range("J7").Value = 0.8238194745364891
$ws.Range("M7").Value = 52.405407
$ws.Range("N7").Value = 157.216221
$ws.Range("O7").Value = 0.9454097874454597
$ws.Range("P7").Value = 0.9454097874454597
$ws.Range("Q7").Value = 12764.24329623794
$ws.Range("R7").Value = 114878.1896661415
$ws.Range("S7").Value = 0.7788469943149725
$ws.Range("T7").Value = 0.7788469943149724
$ws.Range("G8").Value = 51.01955666666666
$ws.Range("H8").Value = 153.05867
$ws.Range("I8").Value = 0.1725638252628821
$ws.Range("J8").Value = 0.1725638252628821
$ws.Range("M8").Value = 1.743137
$ws.Range("N8").Value = 5.229411
$ws.Range("O8").Value = 0.03144673183548247
$ws.Range("P8").Value = 0.03144673183548247
$ws.Range("Q8").Value = 88.93407694926331
$ws.Range("R8").Value = 800.4066925433699
$ws.Range("S8").Value = 0.005426568337546908
$ws.Range("T8").Value = 0.005426568337546908
$ws.Range("G9").Value = 51.01955666666666
$ws.Range("H9").Value = 153.05867
$ws.Range("I9").Value = 0.1725638252628821
$ws.Range("J9").Value = 0.1725638252628821
$ws.Range("N9").Value = 3.848628
$ws.Range("O9").Value = 0.02314348071905789
$ws.Range("P9").Value = 0.02314348071905789
$ws.Range("Q9").Value = 65.45176477830667
$ws.Range("R9").Value = 589.0658830047599
$ws.Range("S9").Value = 0.003993727562778387
$ws.Range("T9").Value = 0.003993727562778386
$ws.Range("G10").Value = 51.01955666666666
$ws.Range("H10").Value = 153.05867
$ws.Range("I10").Value = 0.1725638252628821
$ws.Range("J10").Value = 0.1725638252628821
$ws.Range("M10").Value = 52.405407
$ws.Range("N10").Value = 157.216221
$ws.Range("O10").Value = 0.9454097874454597
$ws.Range("P10").Value = 0.9454097874454597
$ws.Range("Q10").Value = 2673.70063207623
$ws.Range("R10").Value = 24063.30568868607
$ws.Range("S10").Value = 0.1631435293625568
$ws.Range("T10").Value = 0.1631435293625568
